# Add a new row (PF/1.0.3) to the meta-sheet, matching columns dev2/sit2/uat2/prod
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
